# Applies the cryptos list refresh described in the commit
# "Updated cryptos list on Sat Apr 15 10:44:43 UTC 2023 with GitHub Actions".
# Only the Price (D) and Volume(1h) (E) columns change, row by row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.505.82"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "'2.108.14"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D5").Value = "'334.21"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").Value = "'0.5243"
$ws.Range("E7").Value = "  -1.71%  "
$ws.Range("D8").Value = "'0.4526"
$ws.Range("E8").Value = "  +3.11%  "
$ws.Range("D9").Value = "'53.55"
$ws.Range("E9").Value = "  +13.79%  "
$ws.Range("D10").Value = "'0.08992"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").Value = "'1.162"
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("D12").Value = "'24.48"
$ws.Range("E12").Value = "  -2.21%  "
$ws.Range("D13").Value = "'2.102.45"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").Value = "'6.777"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "'96.74"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "'0.06620"
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("D20").Value = "'19.32"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "'6.306"
$ws.Range("E22").Value = "  -0.61%  "
$ws.Range("D23").Value = "'30.557.98"
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("D24").Value = "'12.39"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "'2.353"
$ws.Range("E25").Value = "  +3.33%  "
$ws.Range("D26").Value = "'2.351.84"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").Value = "'22.39"
$ws.Range("D28").Value = "'2.578"
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("D29").Value = "'163.26"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "'132.81"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("E31").Value = "  +1.60%  "
$ws.Range("D32").Value = "'0.1074"
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("D33").Value = "'1.659"
$ws.Range("E33").Value = "  +5.88%  "
$ws.Range("D34").Value = "'6.176"
$ws.Range("D35").Value = "'3.944"
$ws.Range("E35").Value = "  -1.79%  "
$ws.Range("D36").Value = "'10.66"
$ws.Range("E36").Value = "  +12.33%  "
$ws.Range("D37").Value = "'0.02579"
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("D38").Value = "'0.06813"
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("D39").Value = "'5.550"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").Value = "'12.79"
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("D41").Value = "'0.2292"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "'0.6927"
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("D43").Value = "'1.256"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("E44").Value = "  +7.20%  "
$ws.Range("D46").Value = "'0.6421"
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("D47").Value = "'14.06"
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("D48").Value = "'3.663"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "'1.249"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("D50").Value = "'1.218"
$ws.Range("E50").Value = "  +4.16%  "
$ws.Range("E51").Value = "  +0.34%  "
